$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the shared string "space " (with trailing space) -> "space" (no trailing space)
# This is the cell in column A, row 8.
$ws.Range("A8").Value = "space"

# Update the active selection to A8 (as saved in the sheet view)
$ws.Range("A8").Select()
